$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.853.19"
$ws.Range("E2").Value = "'  -0.98%  "

$ws.Range("D3").Value = "'1.871.81"
$ws.Range("E3").Value = "'  -1.35%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.19%  "

$ws.Range("D5").Value = "'301.36"
$ws.Range("E5").Value = "'  -1.81%  "

$ws.Range("E6").Value = "'  -0.16%  "

$ws.Range("D7").Value = "'0.5333"
$ws.Range("E7").Value = "'  +1.91%  "

$ws.Range("D8").Value = "'0.3749"
$ws.Range("E8").Value = "'  -1.48%  "

$ws.Range("D9").Value = "'0.07175"
$ws.Range("E9").Value = "'  -1.48%  "

$ws.Range("D10").Value = "'21.47"
$ws.Range("E10").Value = "'  +0.55%  "

$ws.Range("D11").Value = "'0.8867"
$ws.Range("E11").Value = "'  -1.74%  "

$ws.Range("D12").Value = "'0.08159"
$ws.Range("E12").Value = "'  -0.06%  "

$ws.Range("D13").Value = "'1.876.90"
$ws.Range("E13").Value = "'  -0.69%  "

$ws.Range("D14").Value = "'93.24"
$ws.Range("E14").Value = "'  -2.11%  "

$ws.Range("D15").Value = "'5.261"
$ws.Range("E15").Value = "'  -1.65%  "

$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "'  -0.18%  "

$ws.Range("E17").Value = "'  +0.01%  "

$ws.Range("D18").Value = "'0.000008530"
$ws.Range("E18").Value = "'  -1.30%  "

$ws.Range("E19").Value = "'  -0.13%  "

$ws.Range("D20").Value = "'26.893.05"
$ws.Range("E20").Value = "'  -0.99%  "

$ws.Range("D21").Value = "'4.968"
$ws.Range("E21").Value = "'  -2.80%  "

$ws.Range("E22").Value = "'  -0.90%  "

$ws.Range("D23").Value = "'6.391"
$ws.Range("E23").Value = "'  -1.08%  "

$ws.Range("D24").Value = "'147.07"
$ws.Range("E24").Value = "'  -1.35%  "

$ws.Range("D25").Value = "'2.257"
$ws.Range("E25").Value = "'  -3.32%  "

$ws.Range("D26").Value = "'1.732"
$ws.Range("E26").Value = "'  -0.59%  "

$ws.Range("E27").Value = "'  -1.00%  "

$ws.Range("E28").Value = "'  -1.18%  "

$ws.Range("D29").Value = "'4.741"
$ws.Range("E29").Value = "'  -1.82%  "

$ws.Range("D30").Value = "'4.582"
$ws.Range("E30").Value = "'  -6.26%  "

$ws.Range("D31").Value = "'0.09125"
$ws.Range("E31").Value = "'  -0.94%  "

$ws.Range("D32").Value = "'0.7982"
$ws.Range("E32").Value = "'  +0.72%  "

$ws.Range("D33").Value = "'0.05003"
$ws.Range("E33").Value = "'  -0.73%  "

$ws.Range("D34").Value = "'2.991"
$ws.Range("E34").Value = "'  +0.68%  "

$ws.Range("D35").Value = "'1.171"
$ws.Range("E35").Value = "'  -3.92%  "

$ws.Range("D36").Value = "'0.6019"
$ws.Range("E36").Value = "'  +6.04%  "

$ws.Range("D37").Value = "'2.586"
$ws.Range("E37").Value = "'  -2.35%  "

$ws.Range("E38").Value = "'  -6.20%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01951"
$ws.Range("E39").Value = "'  -2.06%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.073"
$ws.Range("E40").Value = "'  -0.75%  "

$ws.Range("D41").Value = "'6.613"
$ws.Range("E41").Value = "'  +0.38%  "

$ws.Range("E42").Value = "'  -1.98%  "

$ws.Range("D43").Value = "'115.70"
$ws.Range("E43").Value = "'  -0.62%  "

$ws.Range("D44").Value = "'0.5137"
$ws.Range("E44").Value = "'  +5.43%  "

$ws.Range("D45").Value = "'0.1497"
$ws.Range("E45").Value = "'  -0.94%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "'  -0.19%  "

$ws.Range("D47").Value = "'9.908"
$ws.Range("E47").Value = "'  -2.03%  "

$ws.Range("E48").Value = "'  -0.41%  "

$ws.Range("D49").Value = "'37.62"
$ws.Range("E49").Value = "'  -1.90%  "

$ws.Range("D50").Value = "'0.06017"
$ws.Range("E50").Value = "'  +1.03%  "

$ws.Range("D51").Value = "'62.13"
$ws.Range("E51").Value = "'  -2.77%  "

